$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Bullet: "The most prestigious publications cover CS's intersection..."
#    Insert new qualifying text right after "publications" and before the
#    existing " cover CS's intersection..." text. "Journal Citation Reports"
#    is italicised.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("The most prestigious publications", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(", as defined by SCImago and Thomson Reuters‘ Journal Citation Reports, mainly ")

$rng = $d.Content
$null = $rng.Find.Execute("Journal Citation Reports", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Italic = 1

# ---------------------------------------------------------------------------
# 2. Bullet: "The writing style is clear and direct. ..."
#    Insert new qualifying text right after "The writing style" and before
#    the existing " clear and direct. ..." text. "MIS Quarterly" and
#    "Computer-aided Civil and Structural Engineering" are italicised.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("The writing style", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter(", as evidenced in journals like MIS Quarterly and Computer-aided Civil and Structural Engineering, is")

$rng = $d.Content
$null = $rng.Find.Execute("MIS Quarterly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Italic = 1

$rng = $d.Content
$null = $rng.Find.Execute("Computer-aided Civil and Structural Engineering", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Italic = 1

Write-Output "done"
